# Bitácora.xlsx — "Add files via upload" commit
# Fixes the title typo, rewrites several description cells, fills in the
# remaining log rows (dates + descriptions) and bumps the sheet zoom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Title: "BOTÁCORA" -> "BITÁCORA" -------------------------------------
$ws.Range("D2").Value = "BITÁCORA  "

# --- Existing description cells that were reworded -----------------------
$ws.Range("G6").Value = "Realización de reunión con la docente para recolección de datos del semillero"
$ws.Range("G8").Value = "Realización de requerimientos del sitio web del semillero "
$ws.Range("G9").Value = "Realizacion de Mock ups  del sitio web del semillero "
$ws.Range("G10").Value = "Realización de requerimientos del sitio web del semillero "

# --- New log rows: dates ---------------------------------------------------
$ws.Range("D11").Value = 43745
$ws.Range("D12").Value = 43745
$ws.Range("D13").Value = 43753
$ws.Range("D14").Value = 43754
$ws.Range("D15").Value = 43759

# --- New log rows: descriptions -------------------------------------------
$ws.Range("G11").Value = "Realización de requerimientos del sitio web del semillero "
$ws.Range("G12").Value = "Realizacion de Mock ups  del sitio web del semillero "
$ws.Range("G13").Value = "Realizacion de Mock ups  del sitio web del semillero "
$ws.Range("G14").Value = "Realización de requerimientos del sitio web del semillero "
$ws.Range("G15").Value = "Realización de requerimientos del sitio web del semillero "

# --- Number formatting: rows 7-15 show the date column as dates ----------
$ws.Range("D7:D15").NumberFormat = "d-mmm-yy"

# --- Sheet zoom -------------------------------------------------------------
$ws.Application.ActiveWindow.Zoom = 115
